# Reply-Letter.docx corrections:
#  1) Response to Reviewer comment about hardware description: reword the
#     parenthetical from "(regarding the current approach, not RVC, it will"
#     to "(if we take into account the current industrial approach, it
#     will".
#  2) Response about security: insert "really " before "expose the system".
#  3) Response about security tools: extend "added to the host" with
#     " system and the framework itself" before the closing period.

$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "(regarding the current approach, not RVC, it will ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(if we take into account the current industrial approach, it will ",
    2)
Write-Host "Change 1 found/replaced: $found1"

# --- Change 2 -------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    " Indeed, security is not the focus of the current paper and such a transparent update may expose the system. However, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Indeed, security is not the focus of the current paper and such a transparent update may really expose the system. However, ",
    2)
Write-Host "Change 2 found/replaced: $found2"

# --- Change 3 -------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "Nonetheless, security tools can be added to the host. That matter was briefly discussed in section V.D.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nonetheless, security tools can be added to the host system and the framework itself. That matter was briefly discussed in section V.D.",
    2)
Write-Host "Change 3 found/replaced: $found3"
